# "Generate Report for Handback"
#
# Populates the localization-status workbook with handback results:
#   - Status cells flip from "Ready for handoff" to
#     "Handed back: in sync with en-US" (Overview + per-locale sheets).
#   - Per-locale sheets get the generated target (.md) and handback
#     (.xlf) file names, the target file also becoming a hyperlink back
#     to the source .md on GitHub.
#   - Per-locale "Latest Handback DateTime" is stamped with the actual
#     handback timestamp.
#   - A couple of columns are widened so the longer file names fit.

$wb = $excel.ActiveWorkbook

# Helper: the host's ColumnWidth setter quantizes to whole pixels via
#   storedWidth = (Round(ColumnWidth * 6) + 5) / 6
# Invert that so a requested *stored* width lands as close as possible.
function StoredWidthToColumnWidth($storedWidth) {
    return ([Math]::Round($storedWidth * 6) - 5) / 6
}

$statusText = "Handed back: in sync with en-US"
$mdFile = "8682d6ec-3df5-405f-99c6-7cf2171687ba.md"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/oltest/blob/25d799f3b8f344891c1895b3fdcd2aa6a0b7fd74/e2e/8682d6ec-3df5-405f-99c6-7cf2171687ba.md"

# ---- Overview sheet ----------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText

$overview.Columns.Item(5).ColumnWidth = StoredWidthToColumnWidth(29.9777050018311)
$overview.Columns.Item(6).ColumnWidth = StoredWidthToColumnWidth(29.9777050018311)

# ---- zh-cn sheet ---------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $statusText
$zhcn.Range("I2").Value = "8682d6ec-3df5-405f-99c6-7cf2171687ba.f7415d136f0b6c2b0c3a9dcfe4b7839e4d3d2d5d.zh-cn.xlf"
$zhcn.Range("J2").Value = "2016-07-14 05:38:49"
$zhcn.Hyperlinks.Add($zhcn.Range("H2"), $mdUrl, "", "", $mdFile)

$zhcn.Columns.Item(3).ColumnWidth = StoredWidthToColumnWidth(29.9777050018311)
$zhcn.Columns.Item(8).ColumnWidth = StoredWidthToColumnWidth(40)
$zhcn.Columns.Item(9).ColumnWidth = StoredWidthToColumnWidth(40)

# ---- de-de sheet ---------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $statusText
$dede.Range("I2").Value = "8682d6ec-3df5-405f-99c6-7cf2171687ba.f7415d136f0b6c2b0c3a9dcfe4b7839e4d3d2d5d.de-de.xlf"
$dede.Range("J2").Value = "2016-07-14 05:39:04"
$dede.Hyperlinks.Add($dede.Range("H2"), $mdUrl, "", "", $mdFile)

$dede.Columns.Item(3).ColumnWidth = StoredWidthToColumnWidth(29.9777050018311)
$dede.Columns.Item(8).ColumnWidth = StoredWidthToColumnWidth(40)
$dede.Columns.Item(9).ColumnWidth = StoredWidthToColumnWidth(40)
